# invalid_code_list.xlsx - "Adjust tests for new restrictions"
#
# The two bogus/placeholder shared strings "CRC001asdas4LM" and
# "CRC0014LM_2" are no longer used anywhere; the model_id values that used
# to hold them are replaced with the well-formed id "CRC0228PRaS" (and the
# sibling "CRC0014LM" model ids become "CRC0228PR"), so both stray strings
# fall out of the shared-string table once nothing references them anymore.
# This touches the "cell_model" and "sharing" sheets.

$wb = $excel.ActiveWorkbook

$wsPatient   = $wb.Worksheets.Item("patient")
$wsCellModel = $wb.Worksheets.Item("cell_model")
$wsSharing   = $wb.Worksheets.Item("sharing")

# --- cell_model sheet -------------------------------------------------
# A2/A3 need the plain "normal" cell style (same as the rest of the grid,
# cellXfs index 4) instead of the heavier bordered style they had before;
# copy formatting from a cell that already has the plain style, then set
# the corrected values.
$wsPatient.Range("A2").Copy() | Out-Null
$wsCellModel.Range("A2:A3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsCellModel.Range("A2").Value = "CRC0228PR"
$wsCellModel.Range("A3").Value = "CRC0228PRaS"

# --- sharing sheet ------------------------------------------------------
$wsSharing.Range("A2").Value = "CRC0228PR"
$wsSharing.Range("A3").Value = "CRC0228PRaS"

# --- view / selection state ---------------------------------------------
# Reset the "sharing" tab's scroll position/selection first ...
$wsSharing.Activate()
$wsSharing.Range("A2").Select()

# ... then land on "cell_model", which becomes the active tab and picks up
# tabSelected="1" (and "patient" loses it).
$wsCellModel.Activate()
$wsCellModel.Range("A2").Select()
